# Update "想去人数" (F column) figures for a handful of events.
# These values are duplicated across the "展览" sheet and the aggregated
# "全部类型" sheet, so both need to be updated in lock-step.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Row -> new value for sheet "展览"
$ws1.Range("F8").Value  = 2153
$ws1.Range("F18").Value = 1605
$ws1.Range("F19").Value = 259
$ws1.Range("F22").Value = 262
$ws1.Range("F24").Value = 12342
$ws1.Range("F25").Value = 12393
$ws1.Range("F26").Value = 920
$ws1.Range("F31").Value = 396
$ws1.Range("F32").Value = 1934
$ws1.Range("F35").Value = 209

# Row -> new value for sheet "全部类型"
$ws4.Range("F3").Value  = 82
$ws4.Range("F9").Value  = 2153
$ws4.Range("F12").Value = 82
$ws4.Range("F23").Value = 1605
$ws4.Range("F24").Value = 259
$ws4.Range("F27").Value = 262
$ws4.Range("F29").Value = 12342
$ws4.Range("F30").Value = 12393
$ws4.Range("F31").Value = 920
$ws4.Range("F36").Value = 396
$ws4.Range("F37").Value = 1934
$ws4.Range("F42").Value = 209
